# Update TPM-derived values in the active worksheet, per the new script run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Omg ligand-expressing row)
$ws.Range("G2").Value = 1.10087
$ws.Range("H2").Value = 3.30261
$ws.Range("I2").Value = 0.1843884439613191
$ws.Range("J2").Value = 0.1843884439613191
$ws.Range("M2").Value = 0.04133166666666666
$ws.Range("Q2").Value = 0.04550079188333333
$ws.Range("R2").Value = 0.40950712695
$ws.Range("S2").Value = 0.1843884439613191
$ws.Range("T2").Value = 0.1843884439613191

# Row 3
$ws.Range("I3").Value = 0.4503925067925547
$ws.Range("J3").Value = 0.4503925067925547
$ws.Range("M3").Value = 0.04133166666666666
$ws.Range("S3").Value = 0.4503925067925547
$ws.Range("T3").Value = 0.4503925067925547

# Row 4
$ws.Range("H4").Value = 6.541494999999999
$ws.Range("I4").Value = 0.3652190492461261
$ws.Range("J4").Value = 0.3652190492461262
$ws.Range("M4").Value = 0.04133166666666666
$ws.Range("Q4").Value = 0.09012363028055553
$ws.Range("R4").Value = 0.8111126725249999
$ws.Range("S4").Value = 0.3652190492461261
$ws.Range("T4").Value = 0.3652190492461262
